$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - column F updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 2959
$ws1.Cells.Item(3, 6).Value = 6415
$ws1.Cells.Item(4, 6).Value = 2524
$ws1.Cells.Item(6, 6).Value = 532
$ws1.Cells.Item(7, 6).Value = 60
$ws1.Cells.Item(12, 6).Value = 7375
$ws1.Cells.Item(13, 6).Value = 339
$ws1.Cells.Item(14, 6).Value = 67
$ws1.Cells.Item(16, 6).Value = 245
$ws1.Cells.Item(19, 6).Value = 9026
$ws1.Cells.Item(22, 6).Value = 267
$ws1.Cells.Item(27, 6).Value = 109
$ws1.Cells.Item(29, 6).Value = 19
$ws1.Cells.Item(30, 6).Value = 65
$ws1.Cells.Item(32, 6).Value = 109
$ws1.Cells.Item(37, 6).Value = 1481
$ws1.Cells.Item(38, 6).Value = 754
$ws1.Cells.Item(39, 6).Value = 3879
$ws1.Cells.Item(40, 6).Value = 203
$ws1.Cells.Item(41, 6).Value = 32
$ws1.Cells.Item(43, 6).Value = 67
$ws1.Cells.Item(44, 6).Value = 21
$ws1.Cells.Item(45, 6).Value = 234
$ws1.Cells.Item(46, 6).Value = 9
$ws1.Cells.Item(47, 6).Value = 54
$ws1.Cells.Item(49, 6).Value = 50

# Sheet "演出" (sheet2) - column F updates
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 49
$ws2.Cells.Item(5, 6).Value = 261
$ws2.Cells.Item(8, 6).Value = 28

# Sheet "全部类型" (sheet4) - column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 49
$ws4.Cells.Item(3, 6).Value = 2959
$ws4.Cells.Item(5, 6).Value = 261
$ws4.Cells.Item(6, 6).Value = 6415
$ws4.Cells.Item(7, 6).Value = 2524
$ws4.Cells.Item(10, 6).Value = 532
$ws4.Cells.Item(11, 6).Value = 60
$ws4.Cells.Item(15, 6).Value = 28
$ws4.Cells.Item(18, 6).Value = 7375
$ws4.Cells.Item(19, 6).Value = 339
$ws4.Cells.Item(20, 6).Value = 67
$ws4.Cells.Item(22, 6).Value = 245
$ws4.Cells.Item(24, 6).Value = 9026
$ws4.Cells.Item(26, 6).Value = 267
$ws4.Cells.Item(30, 6).Value = 109
$ws4.Cells.Item(31, 6).Value = 66
$ws4.Cells.Item(33, 6).Value = 109
$ws4.Cells.Item(38, 6).Value = 1481
$ws4.Cells.Item(39, 6).Value = 754
$ws4.Cells.Item(41, 6).Value = 3879
$ws4.Cells.Item(42, 6).Value = 203
$ws4.Cells.Item(43, 6).Value = 32
$ws4.Cells.Item(46, 6).Value = 234
$ws4.Cells.Item(47, 6).Value = 54
$ws4.Cells.Item(49, 6).Value = 50
